# workflowr::wflow_publish(c("./docs/*", "./analysis/*", "./code/*", "./output/*"))
# Inserts 6 new summary rows ("Control_spatial_binary:" / "Stroke_spatial_binary:"
# blocks) above the existing "Pseudotime" row of the descriptives table, pushing
# the remainder of the table down by 6 rows (old row 24 -> new row 30, ...,
# old row 41 -> new row 47) and growing the sheet's used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24:41 down to 30:47 by inserting 6 blank rows at row 24.
$ws.Rows("24:29").Insert()

# New content for the 6 newly-inserted rows (A..I, row-major).
# Note: "  0.080  " in I24 looks numeric to Excel's auto-typing, so it is
# written with a leading apostrophe (classic "force text" entry) to keep it
# a text value, matching every other p-value cell in the table.
$newRows = @(
    ,@("Control_spatial_binary:", "               ", "               ", "               ", "               ", "              ", "               ", "              ", "'  0.080  ")
    ,@("    FALSE", "  338 (91.1%)  ", "  75 (87.2%)   ", "   57 (100%)   ", "   8 (100%)    ", "   4 (100%)   ", "   15 (100%)   ", "   8 (100%)   ", "         ")
    ,@("    TRUE", "  33 (8.89%)   ", "  11 (12.8%)   ", "   0 (0.00%)   ", "   0 (0.00%)   ", "  0 (0.00%)   ", "   0 (0.00%)   ", "  0 (0.00%)   ", "         ")
    ,@("Stroke_spatial_binary:", "               ", "               ", "               ", "               ", "              ", "               ", "              ", "    .    ")
    ,@("    FALSE", "  240 (64.7%)  ", "  55 (64.0%)   ", "  53 (93.0%)   ", "   8 (100%)    ", "   4 (100%)   ", "  11 (73.3%)   ", "   8 (100%)   ", "         ")
    ,@("    TRUE", "  131 (35.3%)  ", "  31 (36.0%)   ", "   4 (7.02%)   ", "   0 (0.00%)   ", "  0 (0.00%)   ", "   4 (26.7%)   ", "  0 (0.00%)   ", "         ")
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowNum = 24 + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

$wb.Save()
